# Auto-generated Excel COM-interop edit script
# Adds dictionary lookup log rows 80-85 to Sheet1,
# mirroring new client-side search activity (ace, she, young, rose).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A80").Value = "ace"
$ws.Range("B80").Value = "@ace /eis/`r`n*  danh từ`r`n- (đánh bài) quân át, quân xì; điểm 1 (trên quân bài hay con súc sắc)`r`n=duece ace+ một con `"hai`" và một con `"một`" (đánh súc sắc)`r`n- phi công xuất sắc (hạ được trên mười máy bay địch); vận động viên xuất sắc; người giỏi nhất (về cái gì...); nhà vô địch`r`n- (thể dục,thể thao) cú giao bóng thắng điểm; điểm thắng giao bóng (quần vợt)`r`n- chút xíu`r`n=to be within an ace of death+ suýt nữa thì chết`r`n!ace in the hole`r`n- (từ Mỹ,nghĩa Mỹ),  (thông tục) quân bài chủ cao nhất dành cho lúc cần đến (đen & bóng)`r`n- người bạn có thể dựa khi gặp khó khăn`r`n!the ace of aces`r`n- phi công ưu tú nhất`r`n- người xuất sắc nhất trong những người xuất sắc`r`n!the ace of trumps`r`n- quân bài chủ cao nhất`r`n!to have an ace up one's sleeve`r`n- giữ kín quân bài chủ cao nhất dành cho lúc cần`r`n!to trump somebody's ace`r`n- cắt quân át của ai bằng bài chủ`r`n- gạt được một đòn ác hiểm của ai`r`n"
$ws.Rows(80).RowHeight = 374.4

$ws.Range("A81").Value = "she"
$ws.Range("B81").Value = "@she /ʃi:/`r`n*  đại từ`r`n- nó, bà ấy, chị ấy, cô ấy...`r`n=she sings beautifully+ chị ấy hát hay`r`n- nó (chỉ tàu, xe... đã được nhân cách hoá), tàu ấy, xe ấy`r`n=she sails tomorrow+ ngày mai chiếc tàu ấy nhổ neo`r`n- người đàn bà, chị`r`n=she of the black hair+ người đàn bà tóc đen, chị tóc đen`r`n*  danh từ`r`n- đàn bà, con gái`r`n=is the child a he or a she?+ đứa bé là con trai hay con gái?`r`n=the not impossible she+ người có thể yêu được`r`n- con cái`r`n=a litter of two shes and a he+ một ổ hai con cái và một con đực`r`n- (trong từ ghép chỉ động vật) cái`r`n=she-goat+ dê cái`r`n=she-ass+ lừa cái`r`n"
$ws.Rows(81).RowHeight = 273.6

$ws.Range("A82").Value = "young"
$ws.Range("B82").Value = "@young /jʌɳ/`r`n*  tính từ`r`n- trẻ, trẻ tuổi, thiếu niên, thanh niên`r`n=a young man+ một thanh niên`r`n=young people+ thanh niên`r`n=his (her) young woman (man)+ người yêu của nó`r`n=a young family+ gia đình có nhiều con nhỏ`r`n=a young person+ người đàn bà lạ trẻ tuổi thuộc tầng lớp dưới (trong ngôn ngữ những người ở của các gia đình tư sản quý tộc Anh)`r`n=the young person+ những người còn non trẻ ngây thơ cần giữ gìn không cho nghe (đọc) những điều tục tĩu`r`n- non`r`n=young tree+ cây non`r`n- (nghĩa bóng) non trẻ, trẻ tuổi`r`n=a young republic+ nước cộng hoà trẻ tuổi`r`n=he is young for his age+ nó còn non so với tuổi, nó trẻ hơn tuổi`r`n- (nghĩa bóng) non nớt, mới mẻ, chưa có kinh nghiệm`r`n=young in mind+ trí óc còn non nớt`r`n=young in bussiness+ chưa có kinh nghiệm kinh doanh`r`n- (nghĩa bóng) còn sớm, còn ở lúc ban đầu, chưa muộn, chưa quá, chưa già`r`n=the night is young yet+ đêm chưa khuya`r`n=young moon+ trăng non`r`n=autumn is still young+ thu hãy còn đang ở lúc đầu mùa`r`n- của tuổi trẻ, của thời thanh niên, của thời niên thiếu, (thuộc) thế hệ trẻ`r`n=young hope+ hy vọng của tuổi trẻ`r`n=in one's young days+ trong thời kỳ thanh xuân, trong lúc tuổi còn trẻ`r`n- (thông tục) con, nhỏ`r`n=young Smith+ thằng Xmít con, cậu Xmít`r`n*  danh từ`r`n- thú con, chim con (mới đẻ)`r`n=with young+ có chửa (thú)`r`n"
$ws.Rows(82).RowHeight = 409.6

$ws.Range("A83").Value = "ace"
$ws.Range("B83").Value = "@ace /eis/`r`n*  danh từ`r`n- (đánh bài) quân át, quân xì; điểm 1 (trên quân bài hay con súc sắc)`r`n=duece ace+ một con `"hai`" và một con `"một`" (đánh súc sắc)`r`n- phi công xuất sắc (hạ được trên mười máy bay địch); vận động viên xuất sắc; người giỏi nhất (về cái gì...); nhà vô địch`r`n- (thể dục,thể thao) cú giao bóng thắng điểm; điểm thắng giao bóng (quần vợt)`r`n- chút xíu`r`n=to be within an ace of death+ suýt nữa thì chết`r`n!ace in the hole`r`n- (từ Mỹ,nghĩa Mỹ),  (thông tục) quân bài chủ cao nhất dành cho lúc cần đến (đen & bóng)`r`n- người bạn có thể dựa khi gặp khó khăn`r`n!the ace of aces`r`n- phi công ưu tú nhất`r`n- người xuất sắc nhất trong những người xuất sắc`r`n!the ace of trumps`r`n- quân bài chủ cao nhất`r`n!to have an ace up one's sleeve`r`n- giữ kín quân bài chủ cao nhất dành cho lúc cần`r`n!to trump somebody's ace`r`n- cắt quân át của ai bằng bài chủ`r`n- gạt được một đòn ác hiểm của ai`r`n"
$ws.Rows(83).RowHeight = 374.4

$ws.Range("A84").Value = "she"
$ws.Range("B84").Value = "@she /ʃi:/`r`n*  đại từ`r`n- nó, bà ấy, chị ấy, cô ấy...`r`n=she sings beautifully+ chị ấy hát hay`r`n- nó (chỉ tàu, xe... đã được nhân cách hoá), tàu ấy, xe ấy`r`n=she sails tomorrow+ ngày mai chiếc tàu ấy nhổ neo`r`n- người đàn bà, chị`r`n=she of the black hair+ người đàn bà tóc đen, chị tóc đen`r`n*  danh từ`r`n- đàn bà, con gái`r`n=is the child a he or a she?+ đứa bé là con trai hay con gái?`r`n=the not impossible she+ người có thể yêu được`r`n- con cái`r`n=a litter of two shes and a he+ một ổ hai con cái và một con đực`r`n- (trong từ ghép chỉ động vật) cái`r`n=she-goat+ dê cái`r`n=she-ass+ lừa cái`r`n"
$ws.Rows(84).RowHeight = 273.6

$ws.Range("A85").Value = "rose"
$ws.Range("B85").Value = "@rose /rouz/`r`n*  danh từ`r`n- hoa hồng; cây hoa hồng`r`n=a climbing rose+ cây hồng leo`r`n=wild rose+ cây tầm xuân`r`n- cô gái đẹp nhất, hoa khôi`r`n=the rose of the town+ cô gái đẹp nhất tỉnh, hoa khôi của tỉnh`r`n- bông hồng năm cánh (quốc huy của nước Anh)`r`n- màu hồng; (số nhiều) nước da hồng hào`r`n=to have roses in one's cheeks+ má đỏ hồng hào`r`n- nơ hoa hồng (đính ở mũ, ở giầy...)`r`n- hương sen (bình tưới)`r`n- (như) rose-diamond`r`n- (như) rose_window`r`n- chân sừng (phần lồi lên ở gốc sừng nai, hươu...)`r`n- (y học) (the rose) bệnh viêm quầng`r`n=a bed of roses x bed to be born under the rose+ đẻ hoang`r`n=blue rose+ `"bông hồng xanh`" (cái không thể nào có được)`r`n!to gather roses (life's rose)`r`n- tìm thú hưởng lạc`r`n!life is not all roses`r`n- đời không phải hoa hồng cả, đời sống không phải lúc nào cũng sung sướng an nhàn`r`n!a path strewn with roses`r`n- cuộc sống đầy lạc thú`r`n!there is no rose without a thorn`r`n- (tục ngữ) không có hoa hồng nào mà không có gai, không có điều gì sướng mà không có cái khổ kèm theo`r`n!under the rose`r`n- bí mật âm thầm, kín đáo, lén lút`r`n*  tính từ`r`n- hồng, màu hồng`r`n*  ngoại động từ`r`n- nhuộm hồng, nhuốm hồng`r`n=the morning sun rosed the eastern horizon+ mặt trời buổi sáng nhuộm hồng chân trời đằng đông`r`n*  thời quá khứ của rise`r`n"
$ws.Rows(85).RowHeight = 409.6

$wb.Save()
